$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 1413  # H15: was 1370.14
$ws.Cells.Item(15, 9).Value = 1413  # I15: was 1370.14
$ws.Cells.Item(15, 11).Value = 4239  # K15: was 4110.42
$ws.Cells.Item(15, 13).Value = -4070  # M15: was -3941.42
# Row 121
$ws.Cells.Item(121, 8).Value = 923.2222  # H121: was 976.125
$ws.Cells.Item(121, 10).Value = 923.2222  # J121: was 976.125
$ws.Cells.Item(121, 12).Value = 2769.6666  # L121: was 2928.375
$ws.Cells.Item(121, 14).Value = -6263.6666  # N121: was -6422.375
# Row 132
$ws.Cells.Item(132, 8).Value = 30608.15  # H132: was 29870.342
$ws.Cells.Item(132, 9).Value = 33484.61  # I132: was 32589.297
$ws.Cells.Item(132, 11).Value = 100453.83  # K132: was 97767.891
$ws.Cells.Item(132, 13).Value = -97923.83  # M132: was -95237.891
# Row 138
$ws.Cells.Item(138, 8).Value = 2591.6052  # H138: was 2617.4268
$ws.Cells.Item(138, 9).Value = 1834.0588  # I138: was 1865.3939
$ws.Cells.Item(138, 10).Value = 3204.8572  # J138: was 3208.3096
$ws.Cells.Item(138, 11).Value = 5502.1764  # K138: was 5596.1817
$ws.Cells.Item(138, 12).Value = 9614.571599999999  # L138: was 9624.9288
$ws.Cells.Item(138, 13).Value = -362.1764000000003  # M138: was -456.1817000000001
$ws.Cells.Item(138, 14).Value = -19894.5716  # N138: was -19904.9288
# Row 139
$ws.Cells.Item(139, 8).Value = 90000  # H139: was 0
$ws.Cells.Item(139, 10).Value = 90000  # J139: was 0
$ws.Cells.Item(139, 12).Value = 90000  # L139: was 0
$ws.Cells.Item(139, 14).Value = -100280  # N139: was NEW

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5624045  # H32: was 5382570
$ws.Cells.Item(32, 9).Value = 6333165  # I32: was 6254000.5
$ws.Cells.Item(32, 10).Value = 21997.8  # J32: was 19920.46
$ws.Cells.Item(32, 11).Value = 6333165  # K32: was 6254000.5
$ws.Cells.Item(32, 12).Value = 21997.8  # L32: was 19920.46
$ws.Cells.Item(32, 13).Value = -6332878  # M32: was -6253713.5
$ws.Cells.Item(32, 14).Value = -22571.8  # N32: was -20494.46
# Row 103
$ws.Cells.Item(103, 8).Value = 20000  # H103: was 0
$ws.Cells.Item(103, 10).Value = 20000  # J103: was 0
$ws.Cells.Item(103, 12).Value = 20000  # L103: was 0
$ws.Cells.Item(103, 14).Value = -22344  # N103: was NEW

$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Cells.Item(2, 8).Value = 72437  # H2: was 79916
$ws.Cells.Item(2, 10).Value = 72437  # J2: was 79916
$ws.Cells.Item(2, 12).Value = 72437  # L2: was 79916
$ws.Cells.Item(2, 14).Value = -72663  # N2: was -80142
# Row 100
$ws.Cells.Item(100, 8).Value = 40643  # H100: was 20649.2
$ws.Cells.Item(100, 10).Value = 40643  # J100: was 20649.2
$ws.Cells.Item(100, 12).Value = 40643  # L100: was 20649.2
$ws.Cells.Item(100, 14).Value = -42807  # N100: was -22813.2
# Row 107
$ws.Cells.Item(107, 8).Value = 1377  # H107: was 1379.0197
$ws.Cells.Item(107, 9).Value = 1406.0217  # I107: was 1408.2609
$ws.Cells.Item(107, 11).Value = 1406.0217  # K107: was 1408.2609
$ws.Cells.Item(107, 13).Value = 513.9783  # M107: was 511.7391
# Row 134
$ws.Cells.Item(134, 8).Value = 18409.582  # H134: was 18439.521
$ws.Cells.Item(134, 9).Value = 22678.098  # I134: was 23101.78
$ws.Cells.Item(134, 10).Value = 4803.6875  # J134: was 4727
$ws.Cells.Item(134, 11).Value = 68034.29400000001  # K134: was 69305.34
$ws.Cells.Item(134, 12).Value = 14411.0625  # L134: was 14181
$ws.Cells.Item(134, 13).Value = -65499.29400000001  # M134: was -66770.34
$ws.Cells.Item(134, 14).Value = -19481.0625  # N134: was -19251

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 8706.516  # H31: was 8912.968999999999
$ws.Cells.Item(31, 9).Value = 2391.4546  # I31: was 2420.6
$ws.Cells.Item(31, 11).Value = 2391.4546  # K31: was 2420.6
$ws.Cells.Item(31, 13).Value = -2096.4546  # M31: was -2125.6
# Row 34
$ws.Cells.Item(34, 8).Value = 8706.516  # H34: was 8912.968999999999
$ws.Cells.Item(34, 9).Value = 2391.4546  # I34: was 2420.6
$ws.Cells.Item(34, 11).Value = 2391.4546  # K34: was 2420.6
$ws.Cells.Item(34, 13).Value = -2189.4546  # M34: was -2218.6
# Row 43
$ws.Cells.Item(43, 8).Value = 75469.28999999999  # H43: was 86380.836
$ws.Cells.Item(43, 10).Value = 75469.28999999999  # J43: was 86380.836
$ws.Cells.Item(43, 12).Value = 75469.28999999999  # L43: was 86380.836
$ws.Cells.Item(43, 14).Value = -75837.28999999999  # N43: was -86748.836
# Row 56
$ws.Cells.Item(56, 8).Value = 0  # H56: was 25000
$ws.Cells.Item(56, 10).Value = 0  # J56: was 25000
$ws.Cells.Item(56, 12).Value = 0  # L56: was 25000
$ws.Cells.Item(56, 14).ClearContents()  # N56: was -26690
# Row 86
$ws.Cells.Item(86, 8).Value = 2231.5  # H86: was 2013.05
$ws.Cells.Item(86, 9).Value = 2291.875  # I86: was 1953.9286
$ws.Cells.Item(86, 11).Value = 2291.875  # K86: was 1953.9286
$ws.Cells.Item(86, 13).Value = -1168.875  # M86: was -830.9286
# Row 89
$ws.Cells.Item(89, 8).Value = 2231.5  # H89: was 2013.05
$ws.Cells.Item(89, 9).Value = 2291.875  # I89: was 1953.9286
$ws.Cells.Item(89, 11).Value = 11459.375  # K89: was 9769.643
$ws.Cells.Item(89, 13).Value = -5843.375  # M89: was -4153.643
# Row 93
$ws.Cells.Item(93, 8).Value = 16499.75  # H93: was 34000
$ws.Cells.Item(93, 9).Value = 16499.75  # I93: was 34000
$ws.Cells.Item(93, 11).Value = 16499.75  # K93: was 34000
$ws.Cells.Item(93, 13).Value = -14627.75  # M93: was -32128
# Row 101
$ws.Cells.Item(101, 8).Value = 75469.28999999999  # H101: was 86380.836
$ws.Cells.Item(101, 10).Value = 75469.28999999999  # J101: was 86380.836
$ws.Cells.Item(101, 12).Value = 75469.28999999999  # L101: was 86380.836
$ws.Cells.Item(101, 14).Value = -81959.28999999999  # N101: was -92870.836
# Row 141
$ws.Cells.Item(141, 8).Value = 227159.47  # H141: was 227160.14
$ws.Cells.Item(141, 10).Value = 227159.47  # J141: was 227160.14
$ws.Cells.Item(141, 12).Value = 227159.47  # L141: was 227160.14
$ws.Cells.Item(141, 14).Value = -237519.47  # N141: was -237520.14

$ws = $wb.Worksheets.Item("CUL")
# Row 24
$ws.Cells.Item(24, 8).Value = 2593.3076  # H24: was 3141.3
$ws.Cells.Item(24, 9).Value = 766.6667  # I24: was 0
$ws.Cells.Item(24, 11).Value = 2300.0001  # K24: was 0
$ws.Cells.Item(24, 13).Value = -2070.0001  # M24: was NEW
# Row 40
$ws.Cells.Item(40, 8).Value = 160.14285  # H40: was 119
$ws.Cells.Item(40, 9).Value = 139.6  # I40: was 126.5
$ws.Cells.Item(40, 10).Value = 211.5  # J40: was 111.5
$ws.Cells.Item(40, 11).Value = 558.4  # K40: was 506
$ws.Cells.Item(40, 12).Value = 846  # L40: was 446
$ws.Cells.Item(40, 13).Value = -489.4  # M40: was -437
$ws.Cells.Item(40, 14).Value = -984  # N40: was -584
# Row 92
$ws.Cells.Item(92, 8).Value = 632.7  # H92: was 998.1429000000001
$ws.Cells.Item(92, 9).Value = 671  # I92: was 1052.7693
$ws.Cells.Item(92, 11).Value = 2013  # K92: was 3158.3079
$ws.Cells.Item(92, 13).Value = -765  # M92: was -1910.3079
# Row 110
$ws.Cells.Item(110, 8).Value = 21514.572  # H110: was 19825.25
$ws.Cells.Item(110, 9).Value = 2427  # I110: was 5213.5
$ws.Cells.Item(110, 11).Value = 7281  # K110: was 15640.5
$ws.Cells.Item(110, 13).Value = -3191  # M110: was -11550.5
# Row 119
$ws.Cells.Item(119, 8).Value = 630  # H119: was 631
$ws.Cells.Item(119, 9).Value = 630  # I119: was 631
$ws.Cells.Item(119, 11).Value = 1890  # K119: was 1893
$ws.Cells.Item(119, 13).Value = 2948  # M119: was 2945
# Row 134
$ws.Cells.Item(134, 8).Value = 1992.6364  # H134: was 2895.8333
$ws.Cells.Item(134, 9).Value = 1114.875  # I134: was 1343.75
$ws.Cells.Item(134, 10).Value = 4333.3335  # J134: was 6000
$ws.Cells.Item(134, 11).Value = 3344.625  # K134: was 4031.25
$ws.Cells.Item(134, 12).Value = 13000.0005  # L134: was 18000
$ws.Cells.Item(134, 13).Value = 1725.375  # M134: was 1038.75
$ws.Cells.Item(134, 14).Value = -23140.0005  # N134: was -28140

$ws = $wb.Worksheets.Item("GSM")
# Row 126
$ws.Cells.Item(126, 8).Value = 5421.385  # H126: was 5576
$ws.Cells.Item(126, 9).Value = 2581.6667  # I126: was 2500
$ws.Cells.Item(126, 10).Value = 7855.4287  # J126: was 7498.5
$ws.Cells.Item(126, 11).Value = 7745.000100000001  # K126: was 7500
$ws.Cells.Item(126, 12).Value = 23566.2861  # L126: was 22495.5
$ws.Cells.Item(126, 13).Value = -5275.000100000001  # M126: was -5030
$ws.Cells.Item(126, 14).Value = -28506.2861  # N126: was -27435.5
# Row 132
$ws.Cells.Item(132, 8).Value = 247979.61  # H132: was 238286.12
$ws.Cells.Item(132, 9).Value = 377938.7  # I132: was 355753.78
$ws.Cells.Item(132, 11).Value = 1133816.1  # K132: was 1067261.34
$ws.Cells.Item(132, 13).Value = -1131286.1  # M132: was -1064731.34

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 5754.4287  # H7: was 5229.7085
$ws.Cells.Item(7, 9).Value = 5623.316  # I7: was 5175.65
$ws.Cells.Item(7, 10).Value = 7000  # J7: was 5500
$ws.Cells.Item(7, 11).Value = 5623.316  # K7: was 5175.65
$ws.Cells.Item(7, 12).Value = 7000  # L7: was 5500
$ws.Cells.Item(7, 13).Value = -5511.316  # M7: was -5063.65
$ws.Cells.Item(7, 14).Value = -7224  # N7: was -5724
# Row 22
$ws.Cells.Item(22, 8).Value = 49786.715  # H22: was 54812.316
$ws.Cells.Item(22, 9).Value = 111975.664  # I22: was 125911.75
$ws.Cells.Item(22, 10).Value = 3145  # J22: was 3103.6365
$ws.Cells.Item(22, 11).Value = 111975.664  # K22: was 125911.75
$ws.Cells.Item(22, 12).Value = 3145  # L22: was 3103.6365
$ws.Cells.Item(22, 13).Value = -111680.664  # M22: was -125616.75
$ws.Cells.Item(22, 14).Value = -3735  # N22: was -3693.6365
# Row 27
$ws.Cells.Item(27, 8).Value = 49786.715  # H27: was 54812.316
$ws.Cells.Item(27, 9).Value = 111975.664  # I27: was 125911.75
$ws.Cells.Item(27, 10).Value = 3145  # J27: was 3103.6365
$ws.Cells.Item(27, 11).Value = 111975.664  # K27: was 125911.75
$ws.Cells.Item(27, 12).Value = 3145  # L27: was 3103.6365
$ws.Cells.Item(27, 13).Value = -111868.664  # M27: was -125804.75
$ws.Cells.Item(27, 14).Value = -3359  # N27: was -3317.6365
# Row 55
$ws.Cells.Item(55, 8).Value = 844.2941  # H55: was 836.17645
$ws.Cells.Item(55, 10).Value = 1557.875  # J55: was 1540.625
$ws.Cells.Item(55, 12).Value = 1557.875  # L55: was 1540.625
$ws.Cells.Item(55, 14).Value = -1903.875  # N55: was -1886.625
# Row 58
$ws.Cells.Item(58, 8).Value = 3759.8  # H58: was 4902
$ws.Cells.Item(58, 9).Value = 3296.5  # I58: was 4500
$ws.Cells.Item(58, 10).Value = 4068.6667  # J58: was 5103
$ws.Cells.Item(58, 11).Value = 3296.5  # K58: was 4500
$ws.Cells.Item(58, 12).Value = 4068.6667  # L58: was 5103
$ws.Cells.Item(58, 13).Value = -3036.5  # M58: was -4240
$ws.Cells.Item(58, 14).Value = -4588.6667  # N58: was -5623
# Row 100
$ws.Cells.Item(100, 8).Value = 7439.5415  # H100: was 7659.0435
$ws.Cells.Item(100, 9).Value = 2761.8  # I100: was 2788.2856
$ws.Cells.Item(100, 11).Value = 2761.8  # K100: was 2788.2856
$ws.Cells.Item(100, 13).Value = -2220.8  # M100: was -2247.2856
# Row 126
$ws.Cells.Item(126, 8).Value = 5754.4287  # H126: was 5229.7085
$ws.Cells.Item(126, 9).Value = 5623.316  # I126: was 5175.65
$ws.Cells.Item(126, 10).Value = 7000  # J126: was 5500
$ws.Cells.Item(126, 11).Value = 16869.948  # K126: was 15526.95
$ws.Cells.Item(126, 12).Value = 21000  # L126: was 16500
$ws.Cells.Item(126, 13).Value = -14399.948  # M126: was -13056.95
$ws.Cells.Item(126, 14).Value = -25940  # N126: was -21440
# Row 127
$ws.Cells.Item(127, 8).Value = 141984  # H127: was 147232.25
$ws.Cells.Item(127, 10).Value = 141984  # J127: was 147232.25
$ws.Cells.Item(127, 12).Value = 141984  # L127: was 147232.25
$ws.Cells.Item(127, 14).Value = -151904  # N127: was -157152.25

$ws = $wb.Worksheets.Item("WVR")
# Row 98
$ws.Cells.Item(98, 8).Value = 70472  # H98: was 70590
$ws.Cells.Item(98, 9).Value = 30000  # I98: was 0
$ws.Cells.Item(98, 10).Value = 80590  # J98: was 70590
$ws.Cells.Item(98, 11).Value = 30000  # K98: was 0
$ws.Cells.Item(98, 12).Value = 80590  # L98: was 70590
$ws.Cells.Item(98, 13).Value = -27005  # M98: was NEW
$ws.Cells.Item(98, 14).Value = -86580  # N98: was -76580
# Row 113
$ws.Cells.Item(113, 8).Value = 794.86206  # H113: was 794.89655
$ws.Cells.Item(113, 9).Value = 609.4737  # I113: was 640.6667
$ws.Cells.Item(113, 10).Value = 1147.1  # J113: was 1047.2727
$ws.Cells.Item(113, 11).Value = 1828.4211  # K113: was 1922.0001
$ws.Cells.Item(113, 12).Value = 3441.3  # L113: was 3141.8181
$ws.Cells.Item(113, 13).Value = 341.5789  # M113: was 247.9999
$ws.Cells.Item(113, 14).Value = -7781.299999999999  # N113: was -7481.8181
# Row 136
$ws.Cells.Item(136, 8).Value = 12705310  # H136: was 13612790
$ws.Cells.Item(136, 9).Value = 16570762  # I136: was 17323952
$ws.Cells.Item(136, 10).Value = 4538.4287  # J136: was 5197.3335
$ws.Cells.Item(136, 11).Value = 49712286  # K136: was 51971856
$ws.Cells.Item(136, 12).Value = 13615.2861  # L136: was 15592.0005
$ws.Cells.Item(136, 13).Value = -49709736  # M136: was -51969306
$ws.Cells.Item(136, 14).Value = -18715.2861  # N136: was -20692.0005
